$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 194, pushing the existing rows 194:265 down to 195:266.
$ws.Rows("194").Insert()

# Populate the newly inserted row 194 with the new record. The values for
# H,I,K,L,M,N,O,P mirror what used to be in row 194 (now shifted to row 195);
# D (Fecha) and J (Volumen) carry the genuinely new data.
$ws.Range("A194").Value2 = 8
$ws.Range("B194").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C194").Value2 = "Coquimbo"
$ws.Range("D194").Value2 = 44795
$ws.Range("E194").Value2 = 4
$ws.Range("F194").Value2 = 100112031
$ws.Range("G194").Value2 = "Poroto verde"
$ws.Range("H194").Value2 = "Magnum"
$ws.Range("I194").Value2 = "Primera"
$ws.Range("J194").Value2 = 480
$ws.Range("K194").Value2 = 34000
$ws.Range("L194").Value2 = 35000
$ws.Range("M194").Value2 = 34500
$ws.Range("N194").Value2 = "$/malla 25 kilos"
$ws.Range("O194").Value2 = "Perú"
$ws.Range("P194").Value2 = 1380
$ws.Range("Q194").Value2 = 25
$ws.Range("R194").Value2 = "Hortaliza"

# Match the date-format style used by the rest of column D.
$ws.Range("D194").NumberFormat = $ws.Range("D195").NumberFormat
